$wb = $excel.ActiveWorkbook
$wb.DeleteNumberFormat("0.00")
